# Aggiornato Avanzamento.xlsx da lbianco via Streamlit
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Foglio1")

# New roster (sorted A-Z by the normalized "SURNAME NAME" form), replacing
# the previous 23-row table with a fresh 22-row export. One technician
# (GIARDINA SALVATORE) dropped out entirely; several others were renamed
# to the new "SURNAME NAME" convention and got updated hour/price figures.
$data = @(
    ,@("AVOLA IVAN",             96, 35.91, 20)
    ,@("CAMPISI CHRISTIAN",      88, 37.43, 20)
    ,@("CANNARELLA GIUSEPPE ",   56, 47,    20)
    ,@("DI BARTOLOMEO DANIELE",  96, 35.3,  20)
    ,@("DI GIORGIO DANILO",     104, 35.52, 20)
    ,@("FUCILE FRANCESCO",       40, 18.64, 20)
    ,@("GANCI CORRADO",          88, 30.3,  20)
    ,@("GOLINO KEVIN",           96, 40.82, 20)
    ,@("GUASTELLA STEFANO",      96, 36.3,  20)
    ,@("LITTERI DANIELE",       112, 33.36, 20)
    ,@("MANCARELLA SALVATORE",   96, 35.49, 20)
    ,@("D'ANDREA MARCO",         32, 27.39, 20)
    ,@("MAZZARELLA CRISTHIAN",  120, 33.36, 20)
    ,@("PONTE SALVATORE",        32, 27.39, 20)
    ,@("QUADARELLA ANDREA",     104, 29.26, 20)
    ,@("ROCCELLA ALBERTO",       88, 36.36, 20)
    ,@("RUSSO GIOVANNI",         80, 38.1,  20)
    ,@("SANTACROCE CARMELO",     96, 36.18, 20)
    ,@("TRECCARICHI GIOVANNI",   24, 37.76, 20)
    ,@("TUSA MAURIZIO",          80, 34.04, 20)
    ,@("VESPERTINO SIMONE",      96, 40.82, 20)
    ,@("VITTORIO FRANCESCO",    104, 35.25, 20)
)

$oldLastRow = 24
$newLastRow = 1 + $data.Count   # 23

# Drop the row that no longer exists in the refreshed export so every
# subsequent row shifts up by one (keeps the sheet at 23 used rows).
$ws.Rows.Item($oldLastRow).EntireRow.Delete()

# Push the refreshed roster + figures into A2:D23. Column E keeps the
# existing shared "=C-(C*D)/100" formula and just recalculates.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Keep the filter range / autofilter / defined name / dimension in sync
# with the now-smaller table. The sheet already had AutoFilterMode = True,
# so toggle it off first -- re-invoking AutoFilter() on an already-filtered
# sheet just clears it (matches real Excel semantics).
$ws.AutoFilterMode = $false
$ws.Range("A1:D$newLastRow").AutoFilter() | Out-Null

# The hidden _xlnm._FilterDatabase name isn't auto-synced by AutoFilter(),
# so repoint it at the new (smaller) range explicitly.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Foglio1!`$A`$1:`$D`$$newLastRow"
    }
}

$ws.Range("E8").Select() | Out-Null

$wb.Save()
